$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.661.87'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '2.397.55'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = "'564.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').Value = "'141.31"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = "'0.536"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').Value = '2.402.76'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('D13').Value = "'0.341"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').Value = "'26.13"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = "'0.0000169"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.834.66'
$ws.Range('E16').Value = '  +0.57%  '
$ws.Range('D17').Value = '60.475.83'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '2.408.35'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = "'8.10"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.47%  '
$ws.Range('D20').Value = "'10.64"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = "'324.65"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('D23').Value = "'6.05"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D25').Value = "'1.84"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.10%  '
$ws.Range('D26').Value = "'64.93"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = "'570.77"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('D28').Value = "'8.07"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.00%  '
$ws.Range('D29').Value = '2.513.60'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '0.0₃0939'
$ws.Range('E30').Value = '  +1.93%  '
$ws.Range('D31').Value = "'8.09"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.05%  '
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('E34').Value = '  -1.88%  '
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('D36').Value = "'1.46"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.38%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = "'0.371"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = "'152.03"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').Value = "'18.32"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').Value = "'5.16"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = "'2.53"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.44%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = "'1.68"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = "'41.63"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('D46').Value = '0.0₆0288'
$ws.Range('E46').Value = '  +3.72%  '
$ws.Range('D47').Value = "'141.92"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('D48').Value = "'3.57"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').Value = "'0.590"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('D50').Value = "'0.0507"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').Value = "'19.37"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.32%  '
